$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (44) down to the new row (45)
# so the new row reuses existing style indices instead of creating new ones.
$ws.Range("A44:V44").Copy()
$ws.Range("A45:V45").PasteSpecial(-4122)

# Populate the new row's values (Indice 44 = the 45th data row)
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "gibraltar"
$ws.Range("C45").Value = "national-league"
$ws.Range("D45").Value = "2023-2024"
$ws.Range("E45").Value = 45279.875
$ws.Range("F45").Value = "Manchester 62"
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = "Magpies"
$ws.Range("I45").Value = 5
$ws.Range("J45").Value = 4.56
$ws.Range("K45").Value = "19/12/2023 11:46"
$ws.Range("L45").Value = 4.41
$ws.Range("M45").Value = "19/12/2023 20:30"
$ws.Range("N45").Value = 4.52
$ws.Range("O45").Value = "19/12/2023 11:46"
$ws.Range("P45").Value = 4.66
$ws.Range("Q45").Value = "19/12/2023 20:30"
$ws.Range("R45").Value = 1.49
$ws.Range("S45").Value = "19/12/2023 11:46"
$ws.Range("T45").Value = 1.51
$ws.Range("U45").Value = "19/12/2023 20:28"
$ws.Range("V45").Value = "https://www.betexplorer.com/football/gibraltar/national-league/manchester-62-magpies/6ouIHICH/"
